$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.860.80'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.26%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.682.95'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.78%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '556.09'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.13%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '157.80'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.94%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.590'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.85%  '
$ws.Range('E9').Value = '  -2.61%  '
$ws.Range('E10').Value = '  -2.41%  '
$ws.Range('E11').Value = '  -2.47%  '
$ws.Range('E12').Value = '  -3.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.157.10'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.53'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '62.774.12'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.17%  '
$ws.Range('E16').Value = '  -1.99%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.683.90'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -1.89%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.84'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -4.35%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.62'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.31%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '345.10'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -2.31%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.22'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -5.05%  '
$ws.Range('E22').Value = '  +0.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.512'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -1.67%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.27'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.33%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.170'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.43%  '
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.17'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0856'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -5.21%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.38'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +4.73%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.24'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.75%  '
$ws.Range('E31').Value = '  -0.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '164.57'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.35%  '
$ws.Range('E33').Value = '  +0.70%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.998'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('E35').Value = '  +0.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '19.49'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -2.66%  '
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '347.21'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.07%  '
$ws.Range('E39').Value = '  -0.74%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.938'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -3.59%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.99'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.40%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '38.34'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.07%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '20.80'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -4.56%  '
$ws.Range('E44').Value = '  -0.56%  '
$ws.Range('E45').Value = '  -3.97%  '
$ws.Range('E47').Value = '  -4.22%  '
$ws.Range('E48').Value = '  -0.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0970'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -2.90%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '128.86'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -3.36%  '
$ws.Range('E51').Value = '  -3.20%  '
